$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.147.20"
$ws.Range("E2").Value = "  -3.30%  "
$ws.Range("D3").Value = "1.644.09"
$ws.Range("E3").Value = "  -3.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.20"
$ws.Range("E5").Value = "  -2.36%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3909"
$ws.Range("E7").Value = "  -1.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3862"
$ws.Range("E8").Value = "  -4.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.004"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.360"
$ws.Range("E10").Value = "  -7.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.71"
$ws.Range("E11").Value = "  -8.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08467"
$ws.Range("E12").Value = "  -3.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.08"
$ws.Range("E13").Value = "  -7.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.155"
$ws.Range("E14").Value = "  -4.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001286"
$ws.Range("E15").Value = "  -4.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.494"
$ws.Range("E16").Value = "  -5.81%  "
$ws.Range("D17").Value = "1.644.58"
$ws.Range("E17").Value = "  -3.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.48"
$ws.Range("E18").Value = "  -1.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06948"
$ws.Range("E19").Value = "  -3.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.98"
$ws.Range("E20").Value = "  +1.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.958"
$ws.Range("E21").Value = "  -5.08%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.74"
$ws.Range("E23").Value = "  -4.45%  "
$ws.Range("D24").Value = "24.153.67"
$ws.Range("E24").Value = "  -3.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.348"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.729"
$ws.Range("E26").Value = "  -7.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.53"
$ws.Range("E27").Value = "  -5.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.935"
$ws.Range("E28").Value = "  +7.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.60"
$ws.Range("E29").Value = "  -2.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "141.77"
$ws.Range("E30").Value = "  -6.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.371"
$ws.Range("E31").Value = "  -14.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.473"
$ws.Range("E32").Value = "  -5.89%  "
$ws.Range("D33").Value = "1.826.64"
$ws.Range("E33").Value = "  -3.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.230"
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("E35").Value = "  -6.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9834"
$ws.Range("E36").Value = "  -5.06%  "
$ws.Range("E37").Value = "  -6.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2712"
$ws.Range("E38").Value = "  -5.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09253"
$ws.Range("E39").Value = "  -3.43%  "
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.06"
$ws.Range("E41").Value = "  -7.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7630"
$ws.Range("E42").Value = "  -7.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.14"
$ws.Range("E43").Value = "  -6.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.98"
$ws.Range("E44").Value = "  -7.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.493"
$ws.Range("E45").Value = "  -7.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6897"
$ws.Range("E46").Value = "  -6.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.097"
$ws.Range("E47").Value = "  -3.52%  "
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08418"
$ws.Range("E49").Value = "  -4.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.19"
$ws.Range("E50").Value = "  -3.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.265"
$ws.Range("E51").Value = "  -9.14%  "
